# Word COM-interop script implementing the diff:
#  - remove the leading "Below is what I as a tester..." paragraph (+2 blank paragraphs)
#  - fix various typos/wording across the document
#  - collapse a duplicate blank paragraph before "3: Version Control"
#  - collapse a duplicate blank paragraph before "Defect: (10 minutes)"

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Text fixes (paragraph count unaffected by these, so order doesn't matter)
# ---------------------------------------------------------------------------

$d.Content.Find.Execute(
    "To automate a workflow to test the website http.//www.suacedemo.cpm",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "To automate a workflow to test the website http.//www.suacedemo.com",
    2) | Out-Null

$d.Content.Find.Execute(
    "The website was initially tested manually to ensure that it is the correct website as specified. Logi names are provied on the website.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The website was initially tested manually to ensure that it is the correct website as specified. Login names are provided on the website.",
    2) | Out-Null

$d.Content.Find.Execute(
    "The frame work that was decided apon was to use Visual Studio nad Xunit. Selenium was imported onto the project to have the code interact with the website.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The framework that was decided upon was to use Visual Studio and Xunit. Selenium was imported into the project to have the code interact with the website.",
    2) | Out-Null

$d.Content.Find.Execute(
    "All defects raised are to be raised as fowwls:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "All defects raised are to be raised as follows:",
    2) | Out-Null

$d.Content.Find.Execute(
    "d: Whic version of the TAS was the defect found in.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "d: Which version of the TAS was the defect found in.",
    2) | Out-Null

$d.Content.Find.Execute(
    "e: Screen shots of how the defect accurred.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "e: Screen shots of how the defect occurred.",
    2) | Out-Null

$d.Content.Find.Execute(
    "f: Steps to repeat the process. This si required to confirm that the defect is in the SUt and not in the TAS. Verification will have to be done to ensure that the TAS is not actually causing the defect.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "f: Steps to repeat the process. This is required to confirm that the defect is in the SUT and not in the TAS. Verification will have to be done to ensure that the TAS is not actually causing the defect.",
    2) | Out-Null

$d.Content.Find.Execute(
    "1: Date defec found: 11/12/2021",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1: Date defect found: 11/12/2021",
    2) | Out-Null

$d.Content.Find.Execute(
    "5: No screen shots",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "5: No screen shots were captured.",
    2) | Out-Null

$d.Content.Find.Execute(
    "6: The username locked_out_user, was part of the test where all the users are able to log into the SwagLabs website.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "6: The username locked_out_user, was part of the test where all the users are able to log into the SauceLab's website.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Structural fixes: delete paragraphs, bottom-to-top so indices stay valid
# ---------------------------------------------------------------------------

# Remove one of the three blank paragraphs sitting between
# "i: Defect number - ..." and the page break / "Defect: (10 minutes)" heading.
$p = $d.Paragraphs(34)
$rng = $d.Range($p.Range.Start, $p.Range.End)
$rng.Delete()

# Remove one of the two blank paragraphs sitting between the "frame work..."
# paragraph and the "3: Version Control" heading.
$p = $d.Paragraphs(15)
$rng = $d.Range($p.Range.Start, $p.Range.End)
$rng.Delete()

# Remove the leading "Below is what I as a tester..." paragraph together with
# the two blank paragraphs that followed it (paragraphs 1-3).
$p1 = $d.Paragraphs(1)
$p3 = $d.Paragraphs(3)
$rng = $d.Range($p1.Range.Start, $p3.Range.End)
$rng.Delete()
